# Add a new leaderboard row for "Galvan Foods" (customer 0008254), inserted
# at row 28 so the existing rows 28-30 shift down to 29-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28 (VALLEY OFFICE PARK), pushing
# the rows below it down by one.
$ws.Range("A28").EntireRow.Insert()

# Match the row height/formatting used by the rest of the table.
$ws.Rows.Item(28).RowHeight = 13.15

# Populate the new row's data.
$ws.Cells.Item(28, 1).Value = "Galvan Foods"
$ws.Cells.Item(28, 2).Value = "Pietrs, Josh"
$ws.Cells.Item(28, 3).Value = "040"
$ws.Cells.Item(28, 5).Value = "0008254"

# Match the saved selection state (A2:A31, active cell A2).
[void]$ws.Range("A2:A31").Select()
